# The commit inserts a new weekly price record for "Espinaca" (Femacal de
# La Calera) as the new row 554, pushing the previously existing rows
# 554-600 down to 555-601 (the sheet grows from A1:R600 to A1:R601).
#
# Insert a fresh row above row 554 - this shifts every row below it down
# by one, which reproduces the "each old row N becomes new row N+1"
# pattern visible throughout the diff (including the former last row,
# which becomes the new row 601).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(554).Insert()

# Populate the newly inserted row 554 with the new record's data.
$ws.Cells.Item(554, 1).Value  = 3
$ws.Cells.Item(554, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(554, 3).Value  = "Coquimbo"
$ws.Cells.Item(554, 4).Value  = 45166
$ws.Cells.Item(554, 5).Value  = 5
$ws.Cells.Item(554, 6).Value  = 100112012
$ws.Cells.Item(554, 7).Value  = "Espinaca"
$ws.Cells.Item(554, 8).Value  = "Sin especificar"
$ws.Cells.Item(554, 9).Value  = "Primera"
$ws.Cells.Item(554, 10).Value = 110
$ws.Cells.Item(554, 11).Value = 4000
$ws.Cells.Item(554, 12).Value = 4000
$ws.Cells.Item(554, 13).Value = 4000
$ws.Cells.Item(554, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(554, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(554, 16).Value = 1333
$ws.Cells.Item(554, 17).Value = 3
$ws.Cells.Item(554, 18).Value = "Hortaliza"
